$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nlgn2"
$ws.Range("C2").Value = "Nrxn2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.986094666666667
$ws.Range("H2").Value = 5.958284
$ws.Range("I2").Value = 0.105440461398774
$ws.Range("J2").Value = 0.105440461398774
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04996866666666667
$ws.Range("N2").Value = 0.149906
$ws.Range("O2").Value = 0.02333871654029131
$ws.Range("P2").Value = 0.02333871654029131
$ws.Range("Q2").Value = 0.09924250236711112
$ws.Range("R2").Value = 0.8931825213040001
$ws.Range("S2").Value = 0.002460845040463513
$ws.Range("T2").Value = 0.002460845040463513

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nlgn2"
$ws.Range("C3").Value = "Nrxn2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.986094666666667
$ws.Range("H3").Value = 5.958284
$ws.Range("I3").Value = 0.105440461398774
$ws.Range("J3").Value = 0.105440461398774
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.049608666666666
$ws.Range("N3").Value = 6.148826
$ws.Range("O3").Value = 0.9573046246952971
$ws.Range("P3").Value = 0.9573046246952971
$ws.Range("Q3").Value = 4.070716841620444
$ws.Range("R3").Value = 36.636451574584
$ws.Range("S3").Value = 0.1009386413270523
$ws.Range("T3").Value = 0.1009386413270523

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nlgn2"
$ws.Range("C4").Value = "Nrxn2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.986094666666667
$ws.Range("H4").Value = 5.958284
$ws.Range("I4").Value = 0.105440461398774
$ws.Range("J4").Value = 0.105440461398774
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.041443
$ws.Range("N4").Value = 0.124329
$ws.Range("O4").Value = 0.01935665876441155
$ws.Range("P4").Value = 0.01935665876441155
$ws.Range("Q4").Value = 0.08230972127066666
$ws.Range("R4").Value = 0.740787491436
$ws.Range("S4").Value = 0.002040975031258176
$ws.Range("T4").Value = 0.002040975031258175

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nlgn2"
$ws.Range("C5").Value = "Nrxn2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.968049999999998
$ws.Range("H5").Value = 29.90414999999999
$ws.Range("I5").Value = 0.5291972275470833
$ws.Range("J5").Value = 0.5291972275470833
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04996866666666667
$ws.Range("N5").Value = 0.149906
$ws.Range("O5").Value = 0.02333871654029131
$ws.Range("P5").Value = 0.02333871654029131
$ws.Range("Q5").Value = 0.4980901677666666
$ws.Range("R5").Value = 4.482811509899999
$ws.Range("S5").Value = 0.01235078408762942
$ws.Range("T5").Value = 0.01235078408762942

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nlgn2"
$ws.Range("C6").Value = "Nrxn2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.968049999999998
$ws.Range("H6").Value = 29.90414999999999
$ws.Range("I6").Value = 0.5291972275470833
$ws.Range("J6").Value = 0.5291972275470833
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.049608666666666
$ws.Range("N6").Value = 6.148826
$ws.Range("O6").Value = 0.9573046246952971
$ws.Range("P6").Value = 0.9573046246952971
$ws.Range("Q6").Value = 20.43060166976666
$ws.Range("R6").Value = 183.8754150278999
$ws.Range("S6").Value = 0.5066029533067523
$ws.Range("T6").Value = 0.5066029533067523

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nlgn2"
$ws.Range("C7").Value = "Nrxn2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.968049999999998
$ws.Range("H7").Value = 29.90414999999999
$ws.Range("I7").Value = 0.5291972275470833
$ws.Range("J7").Value = 0.5291972275470833
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.041443
$ws.Range("N7").Value = 0.124329
$ws.Range("O7").Value = 0.01935665876441155
$ws.Range("P7").Value = 0.01935665876441155
$ws.Range("Q7").Value = 0.4131058961499999
$ws.Range("R7").Value = 3.717953065349999
$ws.Range("S7").Value = 0.01024349015270154
$ws.Range("T7").Value = 0.01024349015270154

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nlgn2"
$ws.Range("C8").Value = "Nrxn2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.882027333333333
$ws.Range("H8").Value = 20.646082
$ws.Range("I8").Value = 0.3653623110541427
$ws.Range("J8").Value = 0.3653623110541427
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04996866666666667
$ws.Range("N8").Value = 0.149906
$ws.Range("O8").Value = 0.02333871654029131
$ws.Range("P8").Value = 0.02333871654029131
$ws.Range("Q8").Value = 0.3438857298102223
$ws.Range("R8").Value = 3.094971568292
$ws.Range("S8").Value = 0.008527087412198379
$ws.Range("T8").Value = 0.008527087412198379

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nlgn2"
$ws.Range("C9").Value = "Nrxn2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.882027333333333
$ws.Range("H9").Value = 20.646082
$ws.Range("I9").Value = 0.3653623110541427
$ws.Range("J9").Value = 0.3653623110541427
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.049608666666666
$ws.Range("N9").Value = 6.148826
$ws.Range("O9").Value = 0.9573046246952971
$ws.Range("P9").Value = 0.9573046246952971
$ws.Range("Q9").Value = 14.10546286663689
$ws.Range("R9").Value = 126.949165799732
$ws.Range("S9").Value = 0.3497630300614925
$ws.Range("T9").Value = 0.3497630300614925

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Nlgn2"
$ws.Range("C10").Value = "Nrxn2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.882027333333333
$ws.Range("H10").Value = 20.646082
$ws.Range("I10").Value = 0.3653623110541427
$ws.Range("J10").Value = 0.3653623110541427
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.041443
$ws.Range("N10").Value = 0.124329
$ws.Range("O10").Value = 0.01935665876441155
$ws.Range("P10").Value = 0.01935665876441155
$ws.Range("Q10").Value = 0.2852118587753333
$ws.Range("R10").Value = 2.566906728978
$ws.Range("S10").Value = 0.007072193580451831
$ws.Range("T10").Value = 0.00707219358045183
